$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "max" column
$ws.Range("Z1").Value = "max"

# Fill Z2:Z92 with the MAX formula over columns B:Y (using the split form
# B:S,S,T,U,V,W,X,Y exactly as authored) for every data row.
for ($r = 2; $r -le 92; $r++) {
    $cell = $ws.Range("Z$r")
    $cell.Formula = "=MAX(Y$r,X$r,W$r,V$r,U$r,T$r,S$r,B$r`:S$r)"
    $cell.NumberFormat = "0.0"
}

# Add the solution note in AA2
$ws.Range("AA2").Value = "ЕСЛИ(B2:Y2 < (Z2/2))"

# Adjust the view: scroll so column P is at the left edge, zoom to 120%,
# and finish with AA2 selected (matches the author's final view state).
$excel.Goto($ws.Range("P1"), $true)
$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 120
$ws.Range("AA2").Select()
